# Rename the worksheet to its full legal entity name.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "SV Health Investors, LLC"

# The now-empty column L (it only ever carried formatting, no data) is
# removed entirely - this shifts the used range back to A1:K5.
$ws.Range("L:L").Delete()

# Strip the leftover per-cell formatting (explicit black font) from the
# remaining data so the sheet parts still work when the workbook/sheet
# is otherwise empty of custom styling.
$ws.Range("A1:K5").ClearFormats()
